$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Center the three title-page paragraphs (main title, subtitle, author).
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(11).Range.ParagraphFormat.Alignment = 1
$d.Paragraphs.Item(12).Range.ParagraphFormat.Alignment = 1
$d.Paragraphs.Item(13).Range.ParagraphFormat.Alignment = 1

# ---------------------------------------------------------------------------
# 2) Reword the "Cílem práce ..." paragraph: swap "program" and "uživateli"
#    around "bude" / "nabízet", splitting the sentence's second run into five
#    runs along the way, same as a human retyping pieces of the sentence in
#    Word would.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(16)
$pStart = $p.Range.Start

# Rewrite "program bude uživateli nabízet hodiny, ..." -> "bude program
# nabízet uživateli hodiny, ..." across the whole tail of the paragraph so it
# collapses into a single freshly-written run.
$full = $p.Range.Text
$idx = $full.IndexOf("program bude uživateli nabízet")
$matchStart = $pStart + $idx
$matchEnd = $p.Range.End - 1
$r = $d.Range($matchStart, $matchEnd)
$r.Bold = $true
$r.Text = "bude program nabízet uživateli hodiny, které je možné do rozvrhu přidat. "
$r.Bold = $false

# Split "bude " off the front.
$full = $p.Range.Text
$idx = $full.IndexOf("bude program")
$s = $pStart + $idx
$e = $s + "bude ".Length
$r = $d.Range($s, $e)
$r.Bold = $true
$r.Bold = $false

# Split "program" off next.
$full = $p.Range.Text
$idx = $full.IndexOf("program nabízet")
$s = $pStart + $idx
$e = $s + "program".Length
$r = $d.Range($s, $e)
$r.Bold = $true
$r.Bold = $false

# Split " nabízet " off next, leaving "uživateli hodiny, ..." as the final run.
$full = $p.Range.Text
$idx = $full.IndexOf(" nabízet uživateli")
$s = $pStart + $idx
$e = $s + " nabízet ".Length
$r = $d.Range($s, $e)
$r.Bold = $true
$r.Bold = $false
